$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A92").Value = 45946
$ws.Range("B92").Value = "四方坪站充电量(kw)"
$ws.Range("C92").Value = 919.90299999999991
$ws.Range("D92").Value = 1034.1699999999998
$ws.Range("E92").Value = 480.23299999999995
$ws.Range("F92").Value = 443.35900000000004
$ws.Range("G92").Value = 353.44100000000003
$ws.Range("H92").Value = 730.66700000000037
$ws.Range("I92").Value = 469.48499999999996
$ws.Range("J92").Value = 228.10999999999999
$ws.Range("K92").Value = 204.62899999999999
$ws.Range("L92").Value = 68.197000000000003
$ws.Range("M92").Value = 142.256
$ws.Range("N92").Value = 121.02000000000001
$ws.Range("O92").Value = 625.00400000000013
$ws.Range("P92").Value = 1303.299
$ws.Range("Q92").Value = 482.36500000000001
$ws.Range("R92").Value = 339.03
$ws.Range("S92").Value = 302.52999999999997
$ws.Range("T92").Value = 218.88799999999998
$ws.Range("U92").Value = 52.32
$ws.Range("V92").Value = 6.18
$ws.Range("W92").Value = 22.44
$ws.Range("X92").Value = 22.62
$ws.Range("Y92").Value = 49.64
$ws.Range("Z92").Value = 17.829000000000001

$ws.Range("A93").Value = 45946
$ws.Range("B93").Value = "高岭站充电量(kw)"
$ws.Range("C93").Value = 470.37599999999992
$ws.Range("D93").Value = 396.26900000000006
$ws.Range("E93").Value = 87.545999999999992
$ws.Range("F93").Value = 80.507000000000005
$ws.Range("G93").Value = 16.68
$ws.Range("H93").Value = 30.046999999999997
$ws.Range("I93").Value = 121.277
$ws.Range("J93").Value = 150.477
$ws.Range("K93").Value = 327.01399999999995
$ws.Range("L93").Value = 233.535
$ws.Range("M93").Value = 134.67400000000001
$ws.Range("N93").Value = 264.495
$ws.Range("O93").Value = 434.62099999999992
$ws.Range("P93").Value = 622.08000000000015
$ws.Range("Q93").Value = 35.010000000000005
$ws.Range("R93").Value = 104.956
$ws.Range("S93").Value = 116.99999999999999
$ws.Range("T93").Value = 125.798
$ws.Range("U93").Value = 35.004000000000005
$ws.Range("V93").Value = 26.949000000000002
$ws.Range("W93").Value = 19.835999999999999
$ws.Range("X93").Value = 101.74499999999999
$ws.Range("Y93").Value = 16.765000000000001
$ws.Range("Z93").Value = 20.077999999999999

$ws.Range("L99").Select()
